$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new row of data (row 33)
$ws.Cells.Item(33, 1).Value = 10002
$ws.Cells.Item(33, 2).Value = 110032
$ws.Cells.Item(33, 3).Value = 10032
$ws.Cells.Item(33, 4).Value = "eng"
$ws.Cells.Item(33, 5).Value = $true
$ws.Cells.Item(33, 6).Value = "superadmin"
$ws.Cells.Item(33, 7).Value = "now()"
$ws.Cells.Item(33, 8).Value = "now()"

# Update the selection to match the target state
$ws.Range("E31").Select()

# Switch calculation mode to manual
$excel.Calculation = -4135
